$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117, shifting row 117 (and below) down to row 118.
$ws.Rows.Item(117).Insert()

# Update row 115 with the new values (previously it held the 44179 record;
# now it holds a brand-new weekly record).
$ws.Cells.Item(115, 4).Value = 44509   # D115 Fecha
$ws.Cells.Item(115, 13).Value = 500    # M115 Volumen
$ws.Cells.Item(115, 14).Value = 9000   # N115 Precio minimo
$ws.Cells.Item(115, 15).Value = 10000  # O115 Precio maximo
$ws.Cells.Item(115, 16).Value = 9500   # P115 Precio promedio ponderado
$ws.Cells.Item(115, 19).Value = 1357   # S115 Precio $/Kg

# Row 116 now holds what used to be row 115's original values (the
# 44179 record with quality "Primera").
$ws.Cells.Item(116, 4).Value = 44179   # D116 Fecha
$ws.Cells.Item(116, 12).Value = "Primera"  # L116 Calidad
$ws.Cells.Item(116, 13).Value = 200    # M116 Volumen
$ws.Cells.Item(116, 14).Value = 11000  # N116 Precio minimo
$ws.Cells.Item(116, 15).Value = 12000  # O116 Precio maximo
$ws.Cells.Item(116, 16).Value = 11500  # P116 Precio promedio ponderado
$ws.Cells.Item(116, 19).Value = 1643   # S116 Precio $/Kg

# New row 117 holds what used to be row 116's values (the 44491 "Especial"
# record).
$ws.Cells.Item(117, 1).Value = 4
$ws.Cells.Item(117, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(117, 3).Value = "Los Lagos"
$ws.Cells.Item(117, 4).Value = 44491
$ws.Cells.Item(117, 5).Value = 10
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100101
$ws.Cells.Item(117, 8).Value = "Berries"
$ws.Cells.Item(117, 9).Value = 100112025
$ws.Cells.Item(117, 10).Value = "Frutilla"
$ws.Cells.Item(117, 11).Value = "Sin especificar"
$ws.Cells.Item(117, 12).Value = "Especial"
$ws.Cells.Item(117, 13).Value = 300
$ws.Cells.Item(117, 14).Value = 12500
$ws.Cells.Item(117, 15).Value = 12500
$ws.Cells.Item(117, 16).Value = 12500
$ws.Cells.Item(117, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(117, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(117, 19).Value = 1786
$ws.Cells.Item(117, 20).Value = 7

# Apply the same date style (style index 2 in original workbook, using the
# custom datetime number format) used by the rest of column D to the new
# row's date cell, mirroring D116/D118.
$ws.Cells.Item(117, 4).NumberFormat = $ws.Cells.Item(118, 4).NumberFormat
